# LeNet9 results.xlsx - fill in CIFAR-10 section (rows 55-77) with measured
# accuracy data (Acck / Acc / RA0 / RAk columns) that had been left blank,
# and correct the Acc0 (D) values for that block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Table 1 (rows 55-62): K=2,R=1,N=3 / K=?,R=2,N=3 groups
# Columns: D=Acc0 E=Acck F=dataset G=LossNum H=Acc I=AccDesc(D-H) J=RA0(H/D) K2=RAk(H/E)
# ---------------------------------------------------------------------

# Row 55 (train, Loss num 0) - D/E corrected, H newly filled, J/K newly added (not shared)
$ws.Range("D55").Value = 0.86250000000000004
$ws.Range("E55").Value = 0.81020000000000003
$ws.Range("H55").Value = 0.81020000000000003
$ws.Range("J55").Formula = "=H55/D55"
$ws.Range("K55").Formula = "=H55/E55"

# Rows 56-58 (same K/R/N group, loss num 1..3) - J/K formulas shared starting row 56
$ws.Range("D56").Value = 0.86250000000000004
$ws.Range("E56").Value = 0.81020000000000003
$ws.Range("H56").Value = 0.6986
$ws.Range("D57").Value = 0.86250000000000004
$ws.Range("E57").Value = 0.81020000000000003
$ws.Range("H57").Value = 0.47649999999999998
$ws.Range("D58").Value = 0.86250000000000004
$ws.Range("E58").Value = 0.81020000000000003
$ws.Range("H58").Value = 0.1
$ws.Range("J56:J77").Formula = "=H56/D56"
$ws.Range("K56:K77").Formula = "=H56/E56"

# Rows 59-62 (test group)
$ws.Range("D59").Value = 0.74319999999999997
$ws.Range("E59").Value = 0.73480000000000001
$ws.Range("H59").Value = 0.73480000000000001
$ws.Range("D60").Value = 0.74319999999999997
$ws.Range("E60").Value = 0.73480000000000001
$ws.Range("H60").Value = 0.65
$ws.Range("D61").Value = 0.74319999999999997
$ws.Range("E61").Value = 0.73480000000000001
$ws.Range("H61").Value = 0.45400000000000001
$ws.Range("D62").Value = 0.74319999999999997
$ws.Range("E62").Value = 0.73480000000000001
$ws.Range("H62").Value = 0.1

# ---------------------------------------------------------------------
# Table 2 (rows 64-77): K=4,R=2,N=6 group
# ---------------------------------------------------------------------

# Row 64 (train, Loss num 0) - individual J/K formulas (first row of block)
$ws.Range("D64").Value = 0.86250000000000004
$ws.Range("E64").Value = 0.85509999999999997
$ws.Range("H64").Value = 0.85509999999999997
$ws.Range("J64").Formula = "=H64/D64"
$ws.Range("K64").Formula = "=H64/E64"

# Rows 65-70 (train group, loss num 1..6)
$ws.Range("D65").Value = 0.86250000000000004
$ws.Range("E65").Value = 0.85509999999999997
$ws.Range("H65").Value = 0.80079999999999996
$ws.Range("D66").Value = 0.86250000000000004
$ws.Range("E66").Value = 0.85509999999999997
$ws.Range("H66").Value = 0.73919999999999997
$ws.Range("D67").Value = 0.86250000000000004
$ws.Range("E67").Value = 0.85509999999999997
$ws.Range("H67").Value = 0.65359999999999996
$ws.Range("D68").Value = 0.86250000000000004
$ws.Range("E68").Value = 0.85509999999999997
$ws.Range("H68").Value = 0.53879999999999995
$ws.Range("D69").Value = 0.86250000000000004
$ws.Range("E69").Value = 0.85509999999999997
$ws.Range("H69").Value = 0.36430000000000001
$ws.Range("D70").Value = 0.86250000000000004
$ws.Range("E70").Value = 0.85509999999999997
$ws.Range("H70").Value = 0.1

# Rows 71-77 (test group, loss num 0..6)
$ws.Range("D71").Value = 0.74319999999999997
$ws.Range("E71").Value = 0.74390000000000001
$ws.Range("H71").Value = 0.74390000000000001
$ws.Range("D72").Value = 0.74319999999999997
$ws.Range("E72").Value = 0.74390000000000001
$ws.Range("H72").Value = 0.71289999999999998
$ws.Range("D73").Value = 0.74319999999999997
$ws.Range("E73").Value = 0.74390000000000001
$ws.Range("H73").Value = 0.66310000000000002

# Row 74 - H74 uses the new percentage/no-fill/side-border style (14 in the
# finished workbook); I/J/K reference the NEXT row's H cell (an off-by-one
# quirk reproduced faithfully from the source data).
$ws.Range("D74").Value = 0.74319999999999997
$ws.Range("E74").Value = 0.74390000000000001
$c = $ws.Range("H74")
$c.NumberFormat = "0.00%"
$c.Borders.Item(7).LineStyle = 1
$c.Borders.Item(10).LineStyle = 1
$c.Borders.Item(8).LineStyle = -4142
$c.Borders.Item(9).LineStyle = -4142
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108
$c.Value = 0.60549999999999993
$ws.Range("I74").Formula = "=D74-H75"
$ws.Range("J74").Formula = "=H75/D74"
$ws.Range("K74").Formula = "=H75/E74"

$ws.Range("D75").Value = 0.74319999999999997
$ws.Range("E75").Value = 0.74390000000000001
$ws.Range("H75").Value = 0.50639999999999996
$ws.Range("I75").Formula = "=D75-H76"
$ws.Range("J75").Formula = "=H76/D75"
$ws.Range("K75").Formula = "=H76/E75"

$ws.Range("D76").Value = 0.74319999999999997
$ws.Range("E76").Value = 0.74390000000000001
$ws.Range("H76").Value = 0.33850000000000002
$ws.Range("I76").Formula = "=D76-H77"
$ws.Range("J76").Formula = "=H77/D76"
$ws.Range("K76").Formula = "=H77/E76"

$ws.Range("D77").Value = 0.74319999999999997
$ws.Range("E77").Value = 0.74390000000000001
$ws.Range("H77").Value = 0.1

# ---------------------------------------------------------------------
# View state: scroll position + active selection (cosmetic, best effort)
# ---------------------------------------------------------------------
$excel.Goto($ws.Range("A61"), $false)
$win = $excel.ActiveWindow
$win.ScrollRow = 61
$win.ScrollColumn = 1
$ws.Range("M77").Select()

Write-Output "LeNet9 CIFAR-10 section updated"
